$wb = $excel.ActiveWorkbook

# --- PIR: rows 404-417 ---
$pirData = New-Object 'object[,]' 14,6
$pirData[0,0] = '2026-01-28'
$pirData[0,1] = '13:01:55'
$pirData[0,2] = '13:00'
$pirData[0,3] = 'Bathroom'
$pirData[0,4] = 'No Motion'
$pirData[0,5] = 'Inactive'
$pirData[1,0] = '2026-01-28'
$pirData[1,1] = '13:01:58'
$pirData[1,2] = '13:00'
$pirData[1,3] = 'Bathroom'
$pirData[1,4] = 'No Motion'
$pirData[1,5] = 'Inactive'
$pirData[2,0] = '2026-01-28'
$pirData[2,1] = '13:01:59'
$pirData[2,2] = '13:00'
$pirData[2,3] = 'Bathroom'
$pirData[2,4] = 'No Motion'
$pirData[2,5] = 'Inactive'
$pirData[3,0] = '2026-01-28'
$pirData[3,1] = '13:02:05'
$pirData[3,2] = '13:00'
$pirData[3,3] = 'Bathroom'
$pirData[3,4] = 'No Motion'
$pirData[3,5] = 'Inactive'
$pirData[4,0] = '2026-01-28'
$pirData[4,1] = '13:02:09'
$pirData[4,2] = '13:00'
$pirData[4,3] = 'Bathroom'
$pirData[4,4] = 'No Motion'
$pirData[4,5] = 'Inactive'
$pirData[5,0] = '2026-01-28'
$pirData[5,1] = '13:02:14'
$pirData[5,2] = '13:00'
$pirData[5,3] = 'Bathroom'
$pirData[5,4] = 'No Motion'
$pirData[5,5] = 'Inactive'
$pirData[6,0] = '2026-01-28'
$pirData[6,1] = '13:02:19'
$pirData[6,2] = '13:00'
$pirData[6,3] = 'Bathroom'
$pirData[6,4] = 'No Motion'
$pirData[6,5] = 'Inactive'
$pirData[7,0] = '2026-01-28'
$pirData[7,1] = '13:02:25'
$pirData[7,2] = '13:00'
$pirData[7,3] = 'Bathroom'
$pirData[7,4] = 'No Motion'
$pirData[7,5] = 'Inactive'
$pirData[8,0] = '2026-01-28'
$pirData[8,1] = '13:02:29'
$pirData[8,2] = '13:00'
$pirData[8,3] = 'Bathroom'
$pirData[8,4] = 'No Motion'
$pirData[8,5] = 'Inactive'
$pirData[9,0] = '2026-01-28'
$pirData[9,1] = '13:02:34'
$pirData[9,2] = '13:00'
$pirData[9,3] = 'Bathroom'
$pirData[9,4] = 'No Motion'
$pirData[9,5] = 'Inactive'
$pirData[10,0] = '2026-01-28'
$pirData[10,1] = '13:02:39'
$pirData[10,2] = '13:00'
$pirData[10,3] = 'Bathroom'
$pirData[10,4] = 'No Motion'
$pirData[10,5] = 'Inactive'
$pirData[11,0] = '2026-01-28'
$pirData[11,1] = '13:02:45'
$pirData[11,2] = '13:00'
$pirData[11,3] = 'Bathroom'
$pirData[11,4] = 'No Motion'
$pirData[11,5] = 'Inactive'
$pirData[12,0] = '2026-01-28'
$pirData[12,1] = '13:02:49'
$pirData[12,2] = '13:00'
$pirData[12,3] = 'Bathroom'
$pirData[12,4] = 'No Motion'
$pirData[12,5] = 'Inactive'
$pirData[13,0] = '2026-01-28'
$pirData[13,1] = '13:02:54'
$pirData[13,2] = '13:00'
$pirData[13,3] = 'Bathroom'
$pirData[13,4] = 'No Motion'
$pirData[13,5] = 'Inactive'
$ws = $wb.Worksheets.Item('PIR')
$rng = $ws.Range("A404:F417")
$rng.NumberFormat = "@"
$rng.Value = $pirData
$rng.Style = "Normal"

# --- Humidity: rows 379-390 ---
$humData = New-Object 'object[,]' 12,6
$humData[0,0] = '2026-01-28'
$humData[0,1] = '13:01:56'
$humData[0,2] = '13:00'
$humData[0,3] = 'Bathroom'
$humData[0,4] = '88.5%'
$humData[0,5] = 'Active'
$humData[1,0] = '2026-01-28'
$humData[1,1] = '13:02:00'
$humData[1,2] = '13:00'
$humData[1,3] = 'Bathroom'
$humData[1,4] = '87.5%'
$humData[1,5] = 'Active'
$humData[2,0] = '2026-01-28'
$humData[2,1] = '13:02:03'
$humData[2,2] = '13:00'
$humData[2,3] = 'Bathroom'
$humData[2,4] = '88.4%'
$humData[2,5] = 'Active'
$humData[3,0] = '2026-01-28'
$humData[3,1] = '13:02:11'
$humData[3,2] = '13:00'
$humData[3,3] = 'Bathroom'
$humData[3,4] = '87.5%'
$humData[3,5] = 'Active'
$humData[4,0] = '2026-01-28'
$humData[4,1] = '13:02:15'
$humData[4,2] = '13:00'
$humData[4,3] = 'Bathroom'
$humData[4,4] = '88.4%'
$humData[4,5] = 'Active'
$humData[5,0] = '2026-01-28'
$humData[5,1] = '13:02:23'
$humData[5,2] = '13:00'
$humData[5,3] = 'Bathroom'
$humData[5,4] = '88.4%'
$humData[5,5] = 'Active'
$humData[6,0] = '2026-01-28'
$humData[6,1] = '13:02:27'
$humData[6,2] = '13:00'
$humData[6,3] = 'Bathroom'
$humData[6,4] = '88.4%'
$humData[6,5] = 'Active'
$humData[7,0] = '2026-01-28'
$humData[7,1] = '13:02:31'
$humData[7,2] = '13:00'
$humData[7,3] = 'Bathroom'
$humData[7,4] = '87.5%'
$humData[7,5] = 'Active'
$humData[8,0] = '2026-01-28'
$humData[8,1] = '13:02:35'
$humData[8,2] = '13:00'
$humData[8,3] = 'Bathroom'
$humData[8,4] = '88.4%'
$humData[8,5] = 'Active'
$humData[9,0] = '2026-01-28'
$humData[9,1] = '13:02:40'
$humData[9,2] = '13:00'
$humData[9,3] = 'Bathroom'
$humData[9,4] = '87.5%'
$humData[9,5] = 'Active'
$humData[10,0] = '2026-01-28'
$humData[10,1] = '13:02:43'
$humData[10,2] = '13:00'
$humData[10,3] = 'Bathroom'
$humData[10,4] = '88.4%'
$humData[10,5] = 'Active'
$humData[11,0] = '2026-01-28'
$humData[11,1] = '13:02:51'
$humData[11,2] = '13:00'
$humData[11,3] = 'Bathroom'
$humData[11,4] = '87.5%'
$humData[11,5] = 'Active'
$ws = $wb.Worksheets.Item('Humidity')
$rng = $ws.Range("A379:F390")
$rng.NumberFormat = "@"
$rng.Value = $humData
$rng.Style = "Normal"

# --- Temperature: rows 379-390 ---
$tempData = New-Object 'object[,]' 12,6
$tempData[0,0] = '2026-01-28'
$tempData[0,1] = '13:01:57'
$tempData[0,2] = '13:00'
$tempData[0,3] = 'Bathroom'
$tempData[0,4] = '22.8C'
$tempData[0,5] = 'Active'
$tempData[1,0] = '2026-01-28'
$tempData[1,1] = '13:02:01'
$tempData[1,2] = '13:00'
$tempData[1,3] = 'Bathroom'
$tempData[1,4] = '22.8C'
$tempData[1,5] = 'Active'
$tempData[2,0] = '2026-01-28'
$tempData[2,1] = '13:02:04'
$tempData[2,2] = '13:00'
$tempData[2,3] = 'Bathroom'
$tempData[2,4] = '22.7C'
$tempData[2,5] = 'Active'
$tempData[3,0] = '2026-01-28'
$tempData[3,1] = '13:02:12'
$tempData[3,2] = '13:00'
$tempData[3,3] = 'Bathroom'
$tempData[3,4] = '22.7C'
$tempData[3,5] = 'Active'
$tempData[4,0] = '2026-01-28'
$tempData[4,1] = '13:02:16'
$tempData[4,2] = '13:00'
$tempData[4,3] = 'Bathroom'
$tempData[4,4] = '22.7C'
$tempData[4,5] = 'Active'
$tempData[5,0] = '2026-01-28'
$tempData[5,1] = '13:02:24'
$tempData[5,2] = '13:00'
$tempData[5,3] = 'Bathroom'
$tempData[5,4] = '22.7C'
$tempData[5,5] = 'Active'
$tempData[6,0] = '2026-01-28'
$tempData[6,1] = '13:02:28'
$tempData[6,2] = '13:00'
$tempData[6,3] = 'Bathroom'
$tempData[6,4] = '22.7C'
$tempData[6,5] = 'Active'
$tempData[7,0] = '2026-01-28'
$tempData[7,1] = '13:02:32'
$tempData[7,2] = '13:00'
$tempData[7,3] = 'Bathroom'
$tempData[7,4] = '22.7C'
$tempData[7,5] = 'Active'
$tempData[8,0] = '2026-01-28'
$tempData[8,1] = '13:02:36'
$tempData[8,2] = '13:00'
$tempData[8,3] = 'Bathroom'
$tempData[8,4] = '22.7C'
$tempData[8,5] = 'Active'
$tempData[9,0] = '2026-01-28'
$tempData[9,1] = '13:02:41'
$tempData[9,2] = '13:00'
$tempData[9,3] = 'Bathroom'
$tempData[9,4] = '22.7C'
$tempData[9,5] = 'Active'
$tempData[10,0] = '2026-01-28'
$tempData[10,1] = '13:02:44'
$tempData[10,2] = '13:00'
$tempData[10,3] = 'Bathroom'
$tempData[10,4] = '22.8C'
$tempData[10,5] = 'Active'
$tempData[11,0] = '2026-01-28'
$tempData[11,1] = '13:02:52'
$tempData[11,2] = '13:00'
$tempData[11,3] = 'Bathroom'
$tempData[11,4] = '22.8C'
$tempData[11,5] = 'Active'
$ws = $wb.Worksheets.Item('Temperature')
$rng = $ws.Range("A379:F390")
$rng.NumberFormat = "@"
$rng.Value = $tempData
$rng.Style = "Normal"

Write-Output "Added PIR rows 404-417, Humidity rows 379-390, Temperature rows 379-390"